$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2009-10")

# --- Row 4 ---
$ws.Range("A4").Value = 1112
$ws.Range("B4").Value = 40190
$ws.Range("C4").Value = "V"
$ws.Range("D4").Value = "Orlando"
$ws.Range("E4").Value = 53
$ws.Range("F4").Value = -2
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 33
$ws.Range("I4").Value = 109
$ws.Range("J4").Value = -180

# --- Row 5 ---
$ws.Range("A5").Value = 1113
$ws.Range("B5").Value = 40190
$ws.Range("C5").Value = "H"
$ws.Range("D5").Value = "Sacramento"
$ws.Range("E5").Value = 58
$ws.Range("F5").Value = -3
$ws.Range("G5").Value = 23
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 88
$ws.Range("J5").Value = 160

# Match existing formatting: column A (bold/border header-like style) and
# column B (date number format) mirror rows 2/3 exactly, so copy their
# formats (not values) down onto the new rows.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B4:B5").PasteSpecial(-4122)

$excel.CutCopyMode = $false
